# Fruta / hortaliza, semanal
# A new weekly record is inserted as row 14 (Macroferia Regional de Talca -
# Arandano (blue), 2021-12-14, "Primera" quality, volume 200), pushing the
# previously existing rows 14-42 down to rows 15-43.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 14; Excel shifts rows 14:42 down to 15:43
# and extends the used range/dimension automatically.
$ws.Rows.Item(14).Insert()

# Populate the new row 14 with the new weekly observation.
$ws.Cells.Item(14, 1).Value  = 5
$ws.Cells.Item(14, 2).Value  = "Macroferia Regional de Talca"
$ws.Cells.Item(14, 3).Value  = "Maule"
$ws.Cells.Item(14, 4).Value  = 44544
$ws.Cells.Item(14, 5).Value  = 7
$ws.Cells.Item(14, 6).Value  = "Fruta"
$ws.Cells.Item(14, 7).Value  = 100101
$ws.Cells.Item(14, 8).Value  = "Berries"
$ws.Cells.Item(14, 9).Value  = 100101001
$ws.Cells.Item(14, 10).Value = "Arándano (blue)"
$ws.Cells.Item(14, 11).Value = "Sin especificar"
$ws.Cells.Item(14, 12).Value = "Primera"
$ws.Cells.Item(14, 13).Value = 200
$ws.Cells.Item(14, 14).Value = 4000
$ws.Cells.Item(14, 15).Value = 4000
$ws.Cells.Item(14, 16).Value = 4000
$ws.Cells.Item(14, 17).Value = "$/bandeja 2 kilos"
$ws.Cells.Item(14, 18).Value = "Provincia de Linares"
$ws.Cells.Item(14, 19).Value = 2000
$ws.Cells.Item(14, 20).Value = 2
